$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "test" column (I) with a default value used when parsing data.
$ws.Range("I1").Value = "test"
$ws.Range("I2").Value = 1

# The "desc" value (JSON blob) that lived in H2 now belongs to row 3 instead.
$ws.Range("H2").ClearContents()
$ws.Range("H3").Value = '{"address": "北京","tag": "常住地"}'

$ws.Range("I4").Value = 3

# Match the selection left behind by the authoring session.
$null = $ws.Range("F17").Select()
